$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("shortage")

$ws.Range("B2").Value = 4.162295600260551
$ws.Range("C2").Value = 1.107274434386906
$ws.Range("D2").Value = 3.759046060306811
$ws.Range("E2").Value = 0.0003162805065090288
$ws.Range("B3").Value = 0.03023121447526297
$ws.Range("C3").Value = 0.08811083674314943
$ws.Range("D3").Value = 0.3431043852572815
$ws.Range("E3").Value = 0.7323874358597695
$ws.Range("B4").Value = 0.05788194520314387
$ws.Range("C4").Value = 0.08423686466845988
$ws.Range("D4").Value = 0.6871331860576255
$ws.Range("E4").Value = 0.49391430881917
$ws.Range("B5").Value = 0.2476551681717035
$ws.Range("C5").Value = 0.0826425806481708
$ws.Range("D5").Value = 2.996701775638284
$ws.Range("E5").Value = 0.003598497022740576
$ws.Range("B6").Value = 0.2576645367456548
$ws.Range("C6").Value = 0.08672326249814835
$ws.Range("D6").Value = 2.971112125205809
$ws.Range("E6").Value = 0.003880380508392888
$ws.Range("B7").Value = -39.44707325299284
$ws.Range("C7").Value = 16.85169152862024
$ws.Range("D7").Value = -2.340837605886478
$ws.Range("E7").Value = 0.02163957395632394
$ws.Range("B8").Value = 110.9833546447109
$ws.Range("C8").Value = 23.33815517252609
$ws.Range("D8").Value = 4.75544677050402
$ws.Range("E8").Value = 0.000008239817215566
$ws.Range("B9").Value = -37.7927646996677
$ws.Range("C9").Value = 25.50455176969453
$ws.Range("D9").Value = -1.481804700625027
$ws.Range("E9").Value = 0.1421789981415316
$ws.Range("B10").Value = -78.8296824093543
$ws.Range("C10").Value = 25.59333074672177
$ws.Range("D10").Value = -3.080086886285855
$ws.Range("E10").Value = 0.002806223942581709
$ws.Range("B11").Value = 60.46218819983819
$ws.Range("C11").Value = 18.6737895074883
$ws.Range("D11").Value = 3.237810310306459
$ws.Range("E11").Value = 0.00173225779679854
$ws.Range("B12").Value = 4.030439022059699
$ws.Range("C12").Value = 0.7190811354684982
$ws.Range("D12").Value = 5.604985061155546
$ws.Range("E12").Value = 0.0000002652903284271935
$ws.Range("B13").Value = -1.882362477267233
$ws.Range("C13").Value = 1.034466677036626
$ws.Range("D13").Value = -1.819645348711979
$ws.Range("E13").Value = 0.07241849061627345
$ws.Range("B14").Value = 0.0696325212998054
$ws.Range("C14").Value = 1.014433579134348
$ws.Range("D14").Value = 0.06864177481114662
$ws.Range("E14").Value = 0.9454398639911507
$ws.Range("B15").Value = -0.9068918658410942
$ws.Range("C15").Value = 0.982265233365745
$ws.Range("D15").Value = -0.9232657687920166
$ws.Range("E15").Value = 0.3585449778807103
$ws.Range("B16").Value = 2.140793443319625
$ws.Range("C16").Value = 0.8037532099556828
$ws.Range("D16").Value = 2.663495979615016
$ws.Range("E16").Value = 0.009287613265408437
